$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Milliman models")
$shp2 = $ws2.Shapes.Item(2)
Write-Output "Name: $($shp2.Name) Width: $($shp2.Width) Left: $($shp2.Left)"
